$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tier 1_obs")

# Clear the library_id_repository values (column I) for data rows 6-34.
# These cells should become empty inline strings (no default mapping).
$ws.Range("I6:I34").ClearContents()

# Append the second SRR run id to library_sequencing_run (column AB) for each row.
$ws.Range("AB6").Value2 = "fastq_generation_SRR13806027, fastq_generation_SRR13806059"
$ws.Range("AB7").Value2 = "fastq_generation_SRR13806028, fastq_generation_SRR13806060"
$ws.Range("AB8").Value2 = "fastq_generation_SRR13806031, fastq_generation_SRR13806063"
$ws.Range("AB9").Value2 = "fastq_generation_SRR13806030, fastq_generation_SRR13806062"
$ws.Range("AB10").Value2 = "fastq_generation_SRR13806029, fastq_generation_SRR13806061"
$ws.Range("AB11").Value2 = "fastq_generation_SRR13806032, fastq_generation_SRR13806064"
$ws.Range("AB12").Value2 = "fastq_generation_SRR13806033, fastq_generation_SRR13806065"
$ws.Range("AB13").Value2 = "fastq_generation_SRR13806037, fastq_generation_SRR13806069"
$ws.Range("AB14").Value2 = "fastq_generation_SRR13806036, fastq_generation_SRR13806067"
$ws.Range("AB15").Value2 = "fastq_generation_SRR13806034, fastq_generation_SRR13806066"
$ws.Range("AB16").Value2 = "fastq_generation_SRR13806038, fastq_generation_SRR13806070"
$ws.Range("AB17").Value2 = "fastq_generation_SRR13806039, fastq_generation_SRR13806071"
$ws.Range("AB18").Value2 = "fastq_generation_SRR13806040, fastq_generation_SRR13806072"
$ws.Range("AB19").Value2 = "fastq_generation_SRR13806043, fastq_generation_SRR13806075"
$ws.Range("AB20").Value2 = "fastq_generation_SRR13806042, fastq_generation_SRR13806074"
$ws.Range("AB21").Value2 = "fastq_generation_SRR13806041, fastq_generation_SRR13806073"
$ws.Range("AB22").Value2 = "fastq_generation_SRR13806044, fastq_generation_SRR13806076"
$ws.Range("AB23").Value2 = "fastq_generation_SRR13806045, fastq_generation_SRR13806077"
$ws.Range("AB24").Value2 = "fastq_generation_SRR13806047, fastq_generation_SRR13806078"
$ws.Range("AB25").Value2 = "fastq_generation_SRR13806023, fastq_generation_SRR13806048"
$ws.Range("AB26").Value2 = "fastq_generation_SRR13806024, fastq_generation_SRR13806049"
$ws.Range("AB27").Value2 = "fastq_generation_SRR13806054, fastq_generation_SRR13806057"
$ws.Range("AB28").Value2 = "fastq_generation_SRR13806058, fastq_generation_SRR13806080"
$ws.Range("AB29").Value2 = "fastq_generation_SRR13806055, fastq_generation_SRR13806068"
$ws.Range("AB30").Value2 = "fastq_generation_SRR13806056, fastq_generation_SRR13806079"
$ws.Range("AB31").Value2 = "fastq_generation_SRR13806025, fastq_generation_SRR13806050"
$ws.Range("AB32").Value2 = "fastq_generation_SRR13806046, fastq_generation_SRR13806053"
$ws.Range("AB33").Value2 = "fastq_generation_SRR13806035, fastq_generation_SRR13806052"
$ws.Range("AB34").Value2 = "fastq_generation_SRR13806026, fastq_generation_SRR13806051"
